$d = $word.ActiveDocument

# Locate the last paragraph in the document body (end of the HUEB-002 user
# story block) -- new content is appended right after it, before sectPr.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)

# Insert a temporary placeholder character so Word creates the new
# paragraph mark, then trim the range back to exclude the paragraph mark
# and clear it -- this yields a truly empty paragraph (no run at all),
# matching the existing blank separator paragraphs used elsewhere in the
# document (e.g. before "HUEB-002").
$r.InsertAfter("`rX")
$blank = $d.Paragraphs.Last
$blankRange = $blank.Range
$blankRange.MoveEnd(1, -1) | Out-Null
$blankRange.Text = ""

# HUEB-003
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertAfter("`rHUEB-003")

# Como cliente del sistema,
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertAfter("`rComo cliente del sistema,")

# para verificar que mi vehículo esté correctamente registrado,
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertAfter("`rpara verificar que mi vehículo esté correctamente registrado,")

# quiero poder consultar el estado actual y los datos asociados de mi vehículo desde mi perfil.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertAfter("`rquiero poder consultar el estado actual y los datos asociados de mi vehículo desde mi perfil.")
